$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 62, shifting the existing row 62 (and everything
# below it, through row 120) down by one to row 63..121.
$ws.Rows(62).Insert()

# Populate the newly inserted row 62 with the new weekly record.
$ws.Range("A62").Value = 8
$ws.Range("B62").Value = "Terminal La Palmera de La Serena"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 44484
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = 100112037
$ws.Range("G62").Value = "Cebollín"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 3060
$ws.Range("K62").Value = 900
$ws.Range("L62").Value = 1000
$ws.Range("M62").Value = 950
$ws.Range("N62").Value = "$/paquete 6 unidades"
$ws.Range("O62").Value = "Provincia del Elquí"
$ws.Range("P62").Value = 158
$ws.Range("Q62").Value = 6
$ws.Range("R62").Value = "Hortaliza"
